$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells L1:N1 - copy the existing header formatting (bold, borders, centered)
# from K1 so the new headers look consistent with the rest of the header row.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# New data values for columns L (apoio_medio), M (contribuicoes), N (media_contribuicoes)
$ws.Range("L2").Value = 90.81853194977892
$ws.Range("M2").Value = 193026
$ws.Range("N2").Value = 329.9589743589743

$ws.Range("L3").Value = 92.63036679831843
$ws.Range("M3").Value = 70527
$ws.Range("N3").Value = 287.865306122449

$ws.Range("L4").Value = 84.63408307975531
$ws.Range("M4").Value = 126119
$ws.Range("N4").Value = 150.6798088410992

$ws.Range("L5").Value = 99.16759340131101
$ws.Range("M5").Value = 77527
$ws.Range("N5").Value = 141.9908424908425

$ws.Range("L6").Value = 18.15137523021585
$ws.Range("M6").Value = 1885
$ws.Range("N6").Value = 14.5

$ws.Range("L7").Value = 27.77589921308953
$ws.Range("M7").Value = 323
$ws.Range("N7").Value = 14.68181818181818
